$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 5384.5
$ws.Cells.Item(5, 9).Value = 3819.2
$ws.Cells.Item(5, 10).Value = 7993.3335
$ws.Cells.Item(5, 11).Value = 3819.2
$ws.Cells.Item(5, 12).Value = 7993.3335
$ws.Cells.Item(5, 13).Value = -3704.2
$ws.Cells.Item(5, 14).Value = -8223.333500000001
$ws.Cells.Item(20, 8).Value = 4505.5
$ws.Cells.Item(20, 9).Value = 1511
$ws.Cells.Item(20, 11).Value = 1511
$ws.Cells.Item(20, 13).Value = -1281
$ws.Cells.Item(21, 8).Value = 4000
$ws.Cells.Item(21, 9).Value = 4000
$ws.Cells.Item(21, 11).Value = 4000
$ws.Cells.Item(21, 13).Value = -3532
$ws.Cells.Item(23, 8).Value = 4000
$ws.Cells.Item(23, 9).Value = 4000
$ws.Cells.Item(23, 11).Value = 4000
$ws.Cells.Item(23, 13).Value = -3766
$ws.Cells.Item(28, 8).Value = 1592.8334
$ws.Cells.Item(28, 9).Value = 524.7646999999999
$ws.Cells.Item(28, 10).Value = 19750
$ws.Cells.Item(28, 11).Value = 524.7646999999999
$ws.Cells.Item(28, 12).Value = 19750
$ws.Cells.Item(28, 13).Value = -39.76469999999995
$ws.Cells.Item(28, 14).Value = -20720
$ws.Cells.Item(32, 8).Value = 15674.75
$ws.Cells.Item(32, 9).Value = 11790
$ws.Cells.Item(32, 10).Value = 18005.6
$ws.Cells.Item(32, 11).Value = 11790
$ws.Cells.Item(32, 12).Value = 18005.6
$ws.Cells.Item(32, 13).Value = -11464
$ws.Cells.Item(32, 14).Value = -18657.6
$ws.Cells.Item(35, 8).Value = 4505.5
$ws.Cells.Item(35, 9).Value = 1511
$ws.Cells.Item(35, 11).Value = 1511
$ws.Cells.Item(35, 13).Value = -1132
$ws.Cells.Item(39, 8).Value = 889.5
$ws.Cells.Item(39, 9).Value = 167.4
$ws.Cells.Item(39, 11).Value = 502.2
$ws.Cells.Item(39, 13).Value = -206.2
$ws.Cells.Item(62, 8).Value = 7961304.5
$ws.Cells.Item(62, 9).Value = 10231678
$ws.Cells.Item(62, 10).Value = 14999.25
$ws.Cells.Item(62, 11).Value = 10231678
$ws.Cells.Item(62, 12).Value = 14999.25
$ws.Cells.Item(62, 13).Value = -10231054
$ws.Cells.Item(62, 14).Value = -16247.25
$ws.Cells.Item(65, 8).Value = 7961304.5
$ws.Cells.Item(65, 9).Value = 10231678
$ws.Cells.Item(65, 10).Value = 14999.25
$ws.Cells.Item(65, 11).Value = 51158390
$ws.Cells.Item(65, 12).Value = 74996.25
$ws.Cells.Item(65, 13).Value = -51155270
$ws.Cells.Item(65, 14).Value = -81236.25
$ws.Cells.Item(74, 8).Value = 7155411.5
$ws.Cells.Item(74, 9).Value = 28576884
$ws.Cells.Item(74, 10).Value = 14920.267
$ws.Cells.Item(74, 11).Value = 28576884
$ws.Cells.Item(74, 12).Value = 14920.267
$ws.Cells.Item(74, 13).Value = -28575948
$ws.Cells.Item(74, 14).Value = -16792.267
$ws.Cells.Item(77, 8).Value = 7155411.5
$ws.Cells.Item(77, 9).Value = 28576884
$ws.Cells.Item(77, 10).Value = 14920.267
$ws.Cells.Item(77, 11).Value = 142884420
$ws.Cells.Item(77, 12).Value = 74601.33499999999
$ws.Cells.Item(77, 13).Value = -142879740
$ws.Cells.Item(77, 14).Value = -83961.33499999999
$ws.Cells.Item(99, 8).Value = 62504060
$ws.Cells.Item(99, 9).Value = 766.53845
$ws.Cells.Item(99, 11).Value = 2299.61535
$ws.Cells.Item(99, 13).Value = -801.61535
$ws.Cells.Item(113, 8).Value = 13514.059
$ws.Cells.Item(113, 9).Value = 5160.2
$ws.Cells.Item(113, 11).Value = 5160.2
$ws.Cells.Item(113, 13).Value = -1906.2
$ws.Cells.Item(116, 8).Value = 6961194.5
$ws.Cells.Item(116, 10).Value = 19197.4
$ws.Cells.Item(116, 12).Value = 19197.4
$ws.Cells.Item(116, 14).Value = -26081.4
$ws.Cells.Item(129, 8).Value = 23827558
$ws.Cells.Item(129, 9).Value = 4377.8
$ws.Cells.Item(129, 10).Value = 37062660
$ws.Cells.Item(129, 11).Value = 13133.4
$ws.Cells.Item(129, 12).Value = 111187980
$ws.Cells.Item(129, 13).Value = -8133.400000000001
$ws.Cells.Item(129, 14).Value = -111197980
$ws.Cells.Item(132, 8).Value = 341856.06
$ws.Cells.Item(132, 9).Value = 403129.1
$ws.Cells.Item(132, 10).Value = 15066.5
$ws.Cells.Item(132, 11).Value = 1209387.3
$ws.Cells.Item(132, 12).Value = 45199.5
$ws.Cells.Item(132, 13).Value = -1206857.3
$ws.Cells.Item(132, 14).Value = -50259.5
$ws.Cells.Item(137, 8).Value = 3092.04
$ws.Cells.Item(137, 9).Value = 2256.2778
$ws.Cells.Item(137, 11).Value = 6768.8334
$ws.Cells.Item(137, 13).Value = -4218.8334

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2098602.8
$ws.Cells.Item(32, 9).Value = 829.02856
$ws.Cells.Item(32, 11).Value = 829.02856
$ws.Cells.Item(32, 13).Value = -542.02856
$ws.Cells.Item(88, 8).Value = 4771.1177
$ws.Cells.Item(88, 9).Value = 1513.25
$ws.Cells.Item(88, 10).Value = 5773.5386
$ws.Cells.Item(88, 11).Value = 1513.25
$ws.Cells.Item(88, 12).Value = 5773.5386
$ws.Cells.Item(88, 13).Value = -1107.25
$ws.Cells.Item(88, 14).Value = -6585.5386
$ws.Cells.Item(91, 8).Value = 4771.1177
$ws.Cells.Item(91, 9).Value = 1513.25
$ws.Cells.Item(91, 10).Value = 5773.5386
$ws.Cells.Item(91, 11).Value = 1513.25
$ws.Cells.Item(91, 12).Value = 5773.5386
$ws.Cells.Item(91, 13).Value = -109.25
$ws.Cells.Item(91, 14).Value = -8581.5386
$ws.Cells.Item(132, 8).Value = 2838931
$ws.Cells.Item(132, 9).Value = 6063658.5
$ws.Cells.Item(132, 10).Value = 151658.33
$ws.Cells.Item(132, 11).Value = 18190975.5
$ws.Cells.Item(132, 12).Value = 454974.99
$ws.Cells.Item(132, 13).Value = -18188445.5
$ws.Cells.Item(132, 14).Value = -460034.99

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 7824.96
$ws.Cells.Item(99, 9).Value = 8154.3613
$ws.Cells.Item(99, 10).Value = 6977.9287
$ws.Cells.Item(99, 11).Value = 8154.3613
$ws.Cells.Item(99, 12).Value = 6977.9287
$ws.Cells.Item(99, 13).Value = -6656.3613
$ws.Cells.Item(99, 14).Value = -9973.9287
$ws.Cells.Item(134, 8).Value = 3405556.5
$ws.Cells.Item(134, 9).Value = 3665234
$ws.Cells.Item(134, 11).Value = 10995702
$ws.Cells.Item(134, 13).Value = -10993167

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2208.9412
$ws.Cells.Item(31, 9).Value = 1174.6666
$ws.Cells.Item(31, 10).Value = 3372.5
$ws.Cells.Item(31, 11).Value = 1174.6666
$ws.Cells.Item(31, 12).Value = 3372.5
$ws.Cells.Item(31, 13).Value = -879.6666
$ws.Cells.Item(31, 14).Value = -3962.5
$ws.Cells.Item(34, 8).Value = 2208.9412
$ws.Cells.Item(34, 9).Value = 1174.6666
$ws.Cells.Item(34, 10).Value = 3372.5
$ws.Cells.Item(34, 11).Value = 1174.6666
$ws.Cells.Item(34, 12).Value = 3372.5
$ws.Cells.Item(34, 13).Value = -972.6666
$ws.Cells.Item(34, 14).Value = -3776.5
$ws.Cells.Item(99, 8).Value = 8551038
$ws.Cells.Item(99, 9).Value = 12349388
$ws.Cells.Item(99, 11).Value = 12349388
$ws.Cells.Item(99, 13).Value = -12347890
$ws.Cells.Item(107, 8).Value = 795.8182
$ws.Cells.Item(107, 9).Value = 615.1111
$ws.Cells.Item(107, 10).Value = 1609
$ws.Cells.Item(107, 11).Value = 615.1111
$ws.Cells.Item(107, 12).Value = 1609
$ws.Cells.Item(107, 13).Value = 1304.8889
$ws.Cells.Item(107, 14).Value = -5449
$ws.Cells.Item(122, 8).Value = 4631
$ws.Cells.Item(122, 9).Value = 4572
$ws.Cells.Item(122, 10).Value = 4660.5
$ws.Cells.Item(122, 11).Value = 13716
$ws.Cells.Item(122, 12).Value = 13981.5
$ws.Cells.Item(122, 13).Value = -11266
$ws.Cells.Item(122, 14).Value = -18881.5
$ws.Cells.Item(126, 8).Value = 8551038
$ws.Cells.Item(126, 9).Value = 12349388
$ws.Cells.Item(126, 11).Value = 37048164
$ws.Cells.Item(126, 13).Value = -37045694
$ws.Cells.Item(132, 8).Value = 7319.2383
$ws.Cells.Item(132, 9).Value = 6379.2144
$ws.Cells.Item(132, 11).Value = 19137.6432
$ws.Cells.Item(132, 13).Value = -16607.6432
$ws.Cells.Item(134, 8).Value = 30311908
$ws.Cells.Item(134, 9).Value = 83336824
$ws.Cells.Item(134, 10).Value = 11954.143
$ws.Cells.Item(134, 11).Value = 250010472
$ws.Cells.Item(134, 12).Value = 35862.429
$ws.Cells.Item(134, 13).Value = -250007937
$ws.Cells.Item(134, 14).Value = -40932.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(26, 8).Value = 359.5625
$ws.Cells.Item(26, 9).Value = 382.57144
$ws.Cells.Item(26, 10).Value = 198.5
$ws.Cells.Item(26, 11).Value = 1147.71432
$ws.Cells.Item(26, 12).Value = 595.5
$ws.Cells.Item(26, 13).Value = -859.71432
$ws.Cells.Item(26, 14).Value = -1171.5
$ws.Cells.Item(81, 8).Value = 2744
$ws.Cells.Item(81, 9).Value = 2283
$ws.Cells.Item(81, 11).Value = 6849
$ws.Cells.Item(81, 13).Value = -5726
$ws.Cells.Item(84, 8).Value = 2744
$ws.Cells.Item(84, 9).Value = 2283
$ws.Cells.Item(84, 11).Value = 20547
$ws.Cells.Item(84, 13).Value = -14931

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 12142.533
$ws.Cells.Item(122, 9).Value = 12499.429
$ws.Cells.Item(122, 11).Value = 37498.287
$ws.Cells.Item(122, 13).Value = -35048.287

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 3999.389
$ws.Cells.Item(40, 9).Value = 3398.5557
$ws.Cells.Item(40, 11).Value = 3398.5557
$ws.Cells.Item(40, 13).Value = -3262.5557
$ws.Cells.Item(93, 8).Value = 1010.8333
$ws.Cells.Item(93, 9).Value = 1285.6364
$ws.Cells.Item(93, 11).Value = 1285.6364
$ws.Cells.Item(93, 13).Value = -37.63640000000009
$ws.Cells.Item(132, 8).Value = 5000
$ws.Cells.Item(132, 9).Value = 5000
$ws.Cells.Item(132, 11).Value = 15000
$ws.Cells.Item(132, 13).Value = -12470
$ws.Cells.Item(133, 8).Value = 87998
$ws.Cells.Item(133, 10).Value = 87998
$ws.Cells.Item(133, 12).Value = 87998
$ws.Cells.Item(133, 14).Value = -93058

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(80, 8).Value = 40000
$ws.Cells.Item(80, 10).Value = 40000
$ws.Cells.Item(80, 12).Value = 40000
$ws.Cells.Item(80, 14).Value = -41996
$ws.Cells.Item(82, 8).Value = 15000
$ws.Cells.Item(82, 9).Value = 15000
$ws.Cells.Item(82, 11).Value = 15000
$ws.Cells.Item(82, 13).Value = -14617
$ws.Cells.Item(83, 8).Value = 40000
$ws.Cells.Item(83, 10).Value = 40000
$ws.Cells.Item(83, 12).Value = 120000
$ws.Cells.Item(83, 14).Value = -129984
$ws.Cells.Item(85, 8).Value = 15000
$ws.Cells.Item(85, 9).Value = 15000
$ws.Cells.Item(85, 11).Value = 15000
$ws.Cells.Item(85, 13).Value = -13674
$ws.Cells.Item(96, 8).Value = 3149.5625
$ws.Cells.Item(96, 10).Value = 3297.1428
$ws.Cells.Item(96, 12).Value = 3297.1428
$ws.Cells.Item(96, 14).Value = -6043.1428
$ws.Cells.Item(119, 8).Value = 34633.332
$ws.Cells.Item(119, 10).Value = 34633.332
$ws.Cells.Item(119, 12).Value = 34633.332
$ws.Cells.Item(119, 14).Value = -44309.332
$ws.Cells.Item(126, 8).Value = 4068.7896
$ws.Cells.Item(126, 9).Value = 2310.889
$ws.Cells.Item(126, 11).Value = 6932.667
$ws.Cells.Item(126, 13).Value = -4462.667
$ws.Cells.Item(132, 8).Value = 5002.3447
$ws.Cells.Item(132, 9).Value = 1829.1765
$ws.Cells.Item(132, 10).Value = 9497.666999999999
$ws.Cells.Item(132, 11).Value = 5487.529500000001
$ws.Cells.Item(132, 12).Value = 28493.001
$ws.Cells.Item(132, 13).Value = -2957.529500000001
$ws.Cells.Item(132, 14).Value = -33553.001
$ws.Cells.Item(136, 8).Value = 17281416
$ws.Cells.Item(136, 9).Value = 21784146
$ws.Cells.Item(136, 11).Value = 65352438
$ws.Cells.Item(136, 13).Value = -65349888
